# Update countries & provincias Spain
# Apply a data refresh + re-sort (by "Casos totales" descending) to the
# Pais sheet, and bump the "updated at" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Rows whose country name and/or stats change position after the refresh.
# Each entry: row, País, Casos totales, Nuevos casos, Casos activos,
#             Recuperados, Casos críticos, Muertes hoy, Muertes
$updates = @(
    @(35,  "Ucrania",                        16847, 422, 4143, 12248, 222, 17, 456),
    @(36,  "Israel",                         16548,   0, 12232, 4052,  61,  0, 264),
    @(67,  "Armenia",                         3860, 142, 1572,  2239,  10,  1,  49),
    @(97,  "El Salvador",                     1112,  75,  405,   687,  19,  0,  20),
    @(98,  "Kirguistan",                      1082,  38,  735,   335,  13,  0,  12),
    @(99,  "Hong Kong",                       1051,   0, 1008,    39,   1,  0,   4),
    @(197, "Mauritania",                        16,   1,    6,     8,   0,  0,   2),
    @(198, "Namibia",                           16,   0,   11,     5,   0,  0,   0),
    @(200, "Curazao",                           16,   0,   14,     1,   0,  0,   1),
    @(201, "Butan",                             15,   4,    5,    10,   0,  0,   0),
    @(215, "Bonaire, San Eustaquio y Saba",       6,   0,    6,     0,   0,  0,   0),
    @(216, "Sahara Occidental",                   6,   0,    6,     0,   0,  0,   0)
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 1).Value = $u[1]
    $ws.Cells.Item($r, 2).Value = $u[2]
    $ws.Cells.Item($r, 3).Value = $u[3]
    $ws.Cells.Item($r, 4).Value = $u[4]
    $ws.Cells.Item($r, 5).Value = $u[5]
    $ws.Cells.Item($r, 6).Value = $u[6]
    $ws.Cells.Item($r, 7).Value = $u[7]
    $ws.Cells.Item($r, 8).Value = $u[8]
}

# Bump the "Datos actualizados" timestamp on A1.
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 09:05"
